# feat: add 2022-Q3 data
#
# 1) "总计" (summary) sheet: insert a new top data row for 2022-Q3 and push
#    the previously existing quarters down by one row.
# 2) Insert a brand-new "2022-Q3" worksheet (before "2022-Q2") holding the
#    per-fund detail rows for that quarter.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Part 1: "总计" sheet - shift rows 2-5 down to 3-6, write new row 2.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

# Shift the B:D (label / count / value) columns down one row at a time,
# working bottom-up so we never overwrite a row before it has been copied.
$summary.Range("B5:D5").Copy($summary.Range("B6:D6"))
$summary.Range("B4:D4").Copy($summary.Range("B5:D5"))
$summary.Range("B3:D3").Copy($summary.Range("B4:D4"))
$summary.Range("B2:D2").Copy($summary.Range("B3:D3"))

# Give the newly-created row 6 the same formatting as the row above it
# (style carries the border/alignment used by column A's index cells).
$summary.Range("A5").Copy($summary.Range("A6"))

# Re-sequence the running index in column A (0,1,2,3,4).
$summary.Range("A2").Value = 0
$summary.Range("A3").Value = 1
$summary.Range("A4").Value = 2
$summary.Range("A5").Value = 3
$summary.Range("A6").Value = 4

# Write the new 2022-Q3 summary values into row 2.
$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 2
$summary.Range("D2").Value = 1.39

# ---------------------------------------------------------------------
# Part 2: new "2022-Q3" worksheet, inserted before "2022-Q2".
# ---------------------------------------------------------------------
$q2Sheet = $wb.Worksheets.Item("2022-Q2")
$q3Sheet = $wb.Worksheets.Add($q2Sheet)
$q3Sheet.Name = "2022-Q3"

# Copy the header row formatting (bold/border style) from the 2022-Q2 sheet
# so the new sheet matches the look of its siblings.
$q2Sheet.Range("B1:H1").Copy($q3Sheet.Range("B1:H1"))
$q2Sheet.Range("A2:H2").Copy($q3Sheet.Range("A2:H2"))
$q2Sheet.Range("A2:H2").Copy($q3Sheet.Range("A3:H3"))

# Header labels.
$q3Sheet.Range("B1").Value = "基金代码"
$q3Sheet.Range("C1").Value = "基金名称"
$q3Sheet.Range("D1").Value = "基金规模"
$q3Sheet.Range("E1").Value = "股票总仓位"
$q3Sheet.Range("F1").Value = "仓位占比"
$q3Sheet.Range("G1").Value = "持有市值(亿元)"
$q3Sheet.Range("H1").Value = "仓位排名"

# Row 2 - 161810 银华内需精选混合（LOF）
$q3Sheet.Range("A2").Value = 0
$q3Sheet.Range("B2").Value = "'161810"
$q3Sheet.Range("C2").Value = "银华内需精选混合（LOF）"
$q3Sheet.Range("D2").Value = "'23.47"
$q3Sheet.Range("E2").Value = "'94.62"
$q3Sheet.Range("F2").Value = "'5.49"
$q3Sheet.Range("G2").Value = "'1.2885"
$q3Sheet.Range("H2").Value = 10

# Row 3 - 180020 银华成长先锋混合
$q3Sheet.Range("A3").Value = 1
$q3Sheet.Range("B3").Value = "'180020"
$q3Sheet.Range("C3").Value = "银华成长先锋混合"
$q3Sheet.Range("D3").Value = "'2.15"
$q3Sheet.Range("E3").Value = "'79.28"
$q3Sheet.Range("F3").Value = "'4.90"
$q3Sheet.Range("G3").Value = "'0.1054"
$q3Sheet.Range("H3").Value = 9
